# Adds a "Copy Activities" sheet between "Create Activities" and "Add Exchanges",
# and updates the headers / sample data on the surrounding sheets so that the
# workbook reflects the new "source_database" / "activity" / "activity_code"
# naming scheme used across the template.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Copy Activities" sheet right after "Create Activities"
#    (i.e. before "Add Exchanges"). NOTE: worksheet references captured
#    before this structural change can become stale/repointed afterwards,
#    so every sheet reference used below is (re-)fetched by name AFTER this
#    insertion happens.
# ---------------------------------------------------------------------------
$wsAddExchangesTmp = $wb.Worksheets.Item("Add Exchanges")
$wsCopy = $wb.Worksheets.Add($wsAddExchangesTmp)
$wsCopy.Name = "Copy Activities"

# Fresh, reliable references for every sheet we touch from here on.
$wsCreate = $wb.Worksheets.Item("Create Activities")
$wsCopy = $wb.Worksheets.Item("Copy Activities")
$wsAddExchanges = $wb.Worksheets.Item("Add Exchanges")
$wsDelete = $wb.Worksheets.Item("Delete Exchanges")

# ---------------------------------------------------------------------------
# 2. "Create Activities" sheet: rename headers, drop the per-row database
#    column (the database is now implied / no longer entered per activity).
# ---------------------------------------------------------------------------
$wsCreate.Range("A1").Value = "activity_database"
$wsCreate.Range("B1").Value = "activity"
$wsCreate.Range("C1").Value = "reference_product"
$wsCreate.Range("D1").Value = "reference_product_amount"
$wsCreate.Range("E1").Value = "reference_product_unit"
$wsCreate.Range("F1").Value = "std_dev"
$wsCreate.Range("G1").Value = "activity_location"
$wsCreate.Range("H1").Value = "activity_version"
$wsCreate.Range("I1").Value = "code"

$wsCreate.Range("A2:A4").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 3. Populate the new "Copy Activities" sheet.
# ---------------------------------------------------------------------------
$wsCopy.Range("A1").Value = "source_database"
$wsCopy.Range("B1").Value = "activity"
$wsCopy.Range("C1").Value = "activity_code"

$wsCopy.Range("A2").Value = "ecoinvent3.7.1 cut-off"
$wsCopy.Range("B2").Value = "borax production, anhydrous, powder"
$wsCopy.Range("C2").Value = "c7bcb2c135dd16e83fd249ad4dc1d966"

$wsCopy.Columns("A:A").ColumnWidth = 20.7109375
$wsCopy.Columns("B:B").ColumnWidth = 35.140625
$wsCopy.Columns("C:C").ColumnWidth = 34.28515625

# ---------------------------------------------------------------------------
# 4. "Add Exchanges" sheet: rename headers, drop the per-row activity
#    database / exchange database columns (only one row keeps the database
#    name, for the cross-database exchange example).
# ---------------------------------------------------------------------------
$wsAddExchanges.Range("A1").Value = "activity_database"
$wsAddExchanges.Range("B1").Value = "exchange_database"
$wsAddExchanges.Range("C1").Value = "activity"
$wsAddExchanges.Range("D1").Value = "activity_code"
$wsAddExchanges.Range("E1").Value = "activity_location"
$wsAddExchanges.Range("F1").Value = "exchange"
$wsAddExchanges.Range("G1").Value = "amount"
$wsAddExchanges.Range("H1").Value = "unit"
$wsAddExchanges.Range("I1").Value = "exchange_location"
$wsAddExchanges.Range("J1").Value = "exchange_type"
$wsAddExchanges.Range("K1").Value = "exchange_code"

$wsAddExchanges.Range("A2:B9").ClearContents() | Out-Null
$wsAddExchanges.Range("B7").Value = "ecoinvent3.7.1 cut-off"

$wsAddExchanges.Columns("D:D").ColumnWidth = 12.7109375
$wsAddExchanges.Columns("J:J").ColumnWidth = 14.5703125

# ---------------------------------------------------------------------------
# 5. "Delete Exchanges" sheet: rename headers (data rows stay the same).
# ---------------------------------------------------------------------------
$wsDelete.Range("A1").Value = "activity_database"
$wsDelete.Range("B1").Value = "activity"
$wsDelete.Range("C1").Value = "activity_code"
$wsDelete.Range("D1").Value = "exchange_database"
$wsDelete.Range("E1").Value = "exchange"
$wsDelete.Range("F1").Value = "exchange_code"

# ---------------------------------------------------------------------------
# 6. Restore per-sheet selections / scroll positions. The "Copy Activities"
#    selection is applied last so that it ends up the active tab, matching
#    the workbook's activeTab staying on index 1 after the insert.
# ---------------------------------------------------------------------------
$wsCreate.Range("A4").Select() | Out-Null
$wsAddExchanges.Range("F7").Select() | Out-Null
$wsDelete.Range("B3").Select() | Out-Null
$wsCopy.Range("C3").Select() | Out-Null
